# Portugal Primeira Liga workbook update (21-04-2024)
# 1) Swap the row-data (columns B:AC, keeping the "id" in column A fixed) for 13
#    pairs of adjacent rows whose match records were stored in the wrong order.
# 2) Insert 3 newly-played matches (21/22 Apr 2024) before the still-unplayed
#    fixtures, shifting those down and renumbering their "id" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

$pairs = @(
    @(14,15), @(48,49), @(97,98), @(102,103), @(123,124), @(167,168),
    @(175,176), @(195,196), @(204,205), @(220,221), @(231,232), @(245,246),
    @(258,259)
)

foreach ($p in $pairs) {
    Swap-Rows $p[0] $p[1]
}

# --- Insert the 3 new match rows at 264..266, pushing the still-unplayed
#     fixtures (old rows 264-266) down to 267-269. ---
$ws.Rows.Item(264).Resize(3).Insert()

function Set-MatchRow($r, $id, $matchId, $date, $home, $away, $fthg, $ftag, $ftr, $vals) {
    $idCell = $ws.Range("A$r")
    $idCell.Font.Bold = $true
    $idCell.HorizontalAlignment = -4108
    $idCell.VerticalAlignment = -4160
    $idCell.Borders.LineStyle = 1
    $idCell.Value2 = $id

    $ws.Range("B$r").Value2 = $matchId
    $ws.Range("C$r").Value2 = "Portugal Primeira Liga"
    $ws.Range("D$r").Value2 = "Portugal Primeira Liga"
    $ws.Range("E$r").Value2 = $date
    $ws.Range("F$r").Value2 = $home
    $ws.Range("G$r").Value2 = $away
    $ws.Range("H$r").Value2 = $fthg
    $ws.Range("I$r").Value2 = $ftag
    $ws.Range("J$r").Value2 = $ftr

    $cols = @("K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value2 = $vals[$i]
    }
}

Set-MatchRow 264 262 6958024 45402.47916666666 "Moreirense" "Gil Vicente" 0 1 "A" @(2.1,3.25,3.75,2.55,3.1,3,0,1.775,2.1,2.25,1.975,1.875,-1,-1,2,-1,1.1,-1,0.875)
Set-MatchRow 265 263 6962299 45402.58333333334 "Boavista" "Estrela" 1 1 "D" @(2.45,3.2,3,2.875,3.1,2.625,0,2.025,1.825,2.25,2.1,1.775,-1,2.1,-1,0,0,-0.5,0.3875)
Set-MatchRow 266 264 6876690 45402.6875 "Braga" "Vizela" 2 1 "H" @(1.3,5.5,9,1.25,6,11,-1.75,1.825,2.025,3.5,1.95,1.9,0.25,-1,-1,-1,1.025,-1,0.8999999999999999)

# The 3 fixtures that were pushed down (old rows 264-266, now 267-269) keep all
# of their data, but their sequential "id" (column A) must advance by 3.
$ws.Range("A267").Value2 = 265
$ws.Range("A268").Value2 = 266
$ws.Range("A269").Value2 = 267

Write-Output "edit complete"
